# Commit: "Added multiple choice answers"
# Insert two new table columns ("answer-2", "answer-3") between the existing
# "answer" and "population" columns of Table1 on Sheet1, populate them with
# quiz answer choices, fix up the "why?" question/answer text for Belgium and
# Switzerland, and refresh the table/column widths/selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1. Physically insert two blank columns at F:G (this shifts the existing
#    population/area/points columns to H:I:J, carrying their styles/values
#    with them).
# ---------------------------------------------------------------------------
$ws.Range("F1:G1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2. Re-grow the table (ListObject) so it spans the two new columns too.
# ---------------------------------------------------------------------------
$lo.Resize($ws.Range("B1:J64"))

# ---------------------------------------------------------------------------
# 3. Header row: make sure every header cell has its correct text (rewriting
#    all of them avoids stale cached ListColumn names).
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "code"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "question"
$ws.Range("E1").Value = "answer"
$ws.Range("F1").Value = "answer-2"
$ws.Range("G1").Value = "answer-3"
$ws.Range("H1").Value = "population"
$ws.Range("I1").Value = "area"
$ws.Range("J1").Value = "points"

# ---------------------------------------------------------------------------
# 4. Fill in the new "answer-2" / "answer-3" choices for every data row.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 64; $r++) {
    $ws.Cells.Item($r, 6).Value = "Eifel tower"
    $ws.Cells.Item($r, 7).Value = "Christ the Redeemer"
}

# ---------------------------------------------------------------------------
# 5. Belgium (row 7) gets a real multiple-choice question/answer about the
#    Manneken Pis statue; Switzerland (row 10) gets the capitalised "Why?"
#    question text.
# ---------------------------------------------------------------------------
$ws.Range("D7").Value = "Whar's the name of a famous statue in Antwerp"
$ws.Range("E7").Value = "Manneken Pis"
$ws.Range("D10").Value = "Why?"

# ---------------------------------------------------------------------------
# 6. Selection / scroll position, matching the author's final cursor spot.
# ---------------------------------------------------------------------------
$ws.Range("F22").Select()

# ---------------------------------------------------------------------------
# 7. Column widths: widen the columns to fit their new contents, like
#    Excel's "AutoFit Column Width" would after adding the new data.
# ---------------------------------------------------------------------------
$ws.Range("C1:C64").EntireColumn.AutoFit()
$ws.Range("D1:D64").EntireColumn.AutoFit()
$ws.Range("E1:E64").EntireColumn.AutoFit()
$ws.Range("G1:G64").EntireColumn.AutoFit()

# ---------------------------------------------------------------------------
# 8. Restore a maximised window (the source workbook had been left minimized).
# ---------------------------------------------------------------------------
$excel.WindowState = -4137
